$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 584.8
$ws.Range("I2").Value = 606.125
$ws.Range("K2").Value = 606.125
$ws.Range("M2").Value = -493.125
$ws.Range("H32").Value = 300
$ws.Range("J32").Value = 300
$ws.Range("L32").Value = 300
$ws.Range("N32").Value = -952
$ws.Range("H33").Value = 558.9
$ws.Range("J33").Value = 998
$ws.Range("L33").Value = 998
$ws.Range("N33").Value = -1456
$ws.Range("H62").Value = 3385.5715
$ws.Range("I62").Value = 3350
$ws.Range("J62").Value = 3433
$ws.Range("K62").Value = 3350
$ws.Range("L62").Value = 3433
$ws.Range("M62").Value = -2726
$ws.Range("N62").Value = -4681
$ws.Range("H65").Value = 3385.5715
$ws.Range("I65").Value = 3350
$ws.Range("J65").Value = 3433
$ws.Range("K65").Value = 16750
$ws.Range("L65").Value = 17165
$ws.Range("M65").Value = -13630
$ws.Range("N65").Value = -23405
$ws.Range("H87").Value = 103354
$ws.Range("J87").Value = 103354
$ws.Range("L87").Value = 103354
$ws.Range("N87").Value = -105850
$ws.Range("H90").Value = 103354
$ws.Range("J90").Value = 103354
$ws.Range("L90").Value = 310062
$ws.Range("N90").Value = -322542
$ws.Range("H94").Value = 10530.25
$ws.Range("I94").Value = 8558
$ws.Range("K94").Value = 8558
$ws.Range("M94").Value = -8107
$ws.Range("H98").Value = 848.3333
$ws.Range("I98").Value = 848.3333
$ws.Range("K98").Value = 848.3333
$ws.Range("M98").Value = 649.6667
$ws.Range("H105").Value = 12750
$ws.Range("J105").Value = 12750
$ws.Range("L105").Value = 12750
$ws.Range("N105").Value = -19738
$ws.Range("H111").Value = 2103.6155
$ws.Range("I111").Value = 3092.5715
$ws.Range("J111").Value = 949.8333
$ws.Range("K111").Value = 9277.7145
$ws.Range("L111").Value = 2849.4999
$ws.Range("M111").Value = -6210.7145
$ws.Range("N111").Value = -8983.499899999999
$ws.Range("H122").Value = 848.3333
$ws.Range("I122").Value = 848.3333
$ws.Range("K122").Value = 2544.9999
$ws.Range("M122").Value = -94.9998999999998
$ws.Range("H132").Value = 42669.75
$ws.Range("I132").Value = 48021.5
$ws.Range("K132").Value = 144064.5
$ws.Range("M132").Value = -141534.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2499.5
$ws.Range("I45").Value = 1999
$ws.Range("K45").Value = 1999
$ws.Range("M45").Value = -1622
$ws.Range("H94").Value = 41333.332
$ws.Range("J94").Value = 41333.332
$ws.Range("L94").Value = 41333.332
$ws.Range("N94").Value = -43135.332
$ws.Range("H132").Value = 3012
$ws.Range("I132").Value = 3012
$ws.Range("K132").Value = 9036
$ws.Range("M132").Value = -6506

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1966.3334
$ws.Range("I105").Value = 1966.3334
$ws.Range("K105").Value = 1966.3334
$ws.Range("M105").Value = -219.3334
$ws.Range("H134").Value = 800
$ws.Range("I134").Value = 600
$ws.Range("K134").Value = 1800
$ws.Range("M134").Value = 735

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 648.5
$ws.Range("I16").Value = 648.5
$ws.Range("K16").Value = 648.5
$ws.Range("M16").Value = -361.5
$ws.Range("H22").Value = 992
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -649
$ws.Range("N22").Value = -1650
$ws.Range("H31").Value = 6980
$ws.Range("J31").Value = 8470.25
$ws.Range("L31").Value = 8470.25
$ws.Range("N31").Value = -9060.25
$ws.Range("H34").Value = 6980
$ws.Range("J34").Value = 8470.25
$ws.Range("L34").Value = 8470.25
$ws.Range("N34").Value = -8874.25
$ws.Range("H105").Value = 1779.5
$ws.Range("I105").Value = 1779.5
$ws.Range("K105").Value = 1779.5
$ws.Range("M105").Value = -32.5
$ws.Range("H113").Value = 648.5
$ws.Range("I113").Value = 648.5
$ws.Range("K113").Value = 648.5
$ws.Range("M113").Value = 1521.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 56946.527
$ws.Range("I4").Value = 1379.5
$ws.Range("J4").Value = 334781.66
$ws.Range("K4").Value = 4138.5
$ws.Range("L4").Value = 1004344.98
$ws.Range("M4").Value = -4026.5
$ws.Range("N4").Value = -1004568.98
$ws.Range("H7").Value = 115.7
$ws.Range("I7").Value = 65.25
$ws.Range("J7").Value = 149.33333
$ws.Range("K7").Value = 195.75
$ws.Range("L7").Value = 447.99999
$ws.Range("M7").Value = -83.75
$ws.Range("N7").Value = -671.99999
$ws.Range("H56").Value = 7000
$ws.Range("I56").Value = 7000
$ws.Range("K56").Value = 7000
$ws.Range("M56").Value = -6470
$ws.Range("H69").Value = 1493
$ws.Range("I69").Value = 1493
$ws.Range("K69").Value = 4479
$ws.Range("M69").Value = -3668
$ws.Range("H72").Value = 1493
$ws.Range("I72").Value = 1493
$ws.Range("K72").Value = 13437
$ws.Range("M72").Value = -9381

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H122").Value = 8499.25
$ws.Range("I122").Value = 8999
$ws.Range("K122").Value = 26997
$ws.Range("M122").Value = -24547

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1644.25
$ws.Range("J55").Value = 1969.25
$ws.Range("L55").Value = 1969.25
$ws.Range("N55").Value = -2315.25
$ws.Range("H61").Value = 1899
$ws.Range("I61").Value = 1899
$ws.Range("K61").Value = 1899
$ws.Range("M61").Value = -1697
$ws.Range("H113").Value = 1899
$ws.Range("I113").Value = 1899
$ws.Range("K113").Value = 1899
$ws.Range("M113").Value = 271
$ws.Range("H122").Value = 4031.0625
$ws.Range("J122").Value = 4563.364
$ws.Range("L122").Value = 13690.092
$ws.Range("N122").Value = -18590.092
$ws.Range("H132").Value = 2504
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 1500
$ws.Range("M132").Value = 1030

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 565.2308
$ws.Range("J100").Value = 483
$ws.Range("L100").Value = 966
$ws.Range("N100").Value = -2048
$ws.Range("H122").Value = 1461.5
$ws.Range("I122").Value = 1199.2
$ws.Range("K122").Value = 3597.6
$ws.Range("M122").Value = -1147.6
$ws.Range("H132").Value = 1112
$ws.Range("I132").Value = 1112
$ws.Range("K132").Value = 3336
$ws.Range("M132").Value = -806
